$d = $word.ActiveDocument

# The paragraph currently holds three separate runs:
#   "<id>"  (Courier New, color 7f6000, sz 18)
#   "p150v_1"  (plain)
#   "</id>"  (Courier New, color 7f6000, sz 18)
# Collapse them into a single run "<id>p150v_1</id>" carrying the
# formatting of the first run, leaving the other "<id>...</id>" pairs
# (fig_p150v_1, fig_p150v_2) untouched.
$d.Content.Find.Execute("<id>p150v_1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p150v_1</id>", 2) | Out-Null
